# Commit: "Fruta / hortaliza, semanal"
#
# The weekly price series for Acelga @ Terminal La Palmera de La Serena
# gets a new (most-recent) observation inserted right before the current
# first data pair (row 160/161), pushing every later row down by two and
# appending two rows at the bottom (duplicating what used to be the last
# pair) -- i.e. a plain row insert that shifts the whole table.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert two blank rows above the existing row 160, shifting 160:291 -> 162:293
$ws.Rows("160:161").Insert()

# New "Primera" quality observation
$ws.Range("A160").Value = 8
$ws.Range("B160").Value = "Terminal La Palmera de La Serena"
$ws.Range("C160").Value = "Coquimbo"
$ws.Range("D160").Value = 44554
$ws.Range("E160").Value = 4
$ws.Range("F160").Value = 100112009
$ws.Range("G160").Value = "Acelga"
$ws.Range("H160").Value = "Sin especificar"
$ws.Range("I160").Value = "Primera"
$ws.Range("J160").Value = 2400
$ws.Range("K160").Value = 450
$ws.Range("L160").Value = 500
$ws.Range("M160").Value = 475
$ws.Range("N160").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O160").Value = "Provincia del Elquí"
$ws.Range("P160").Value = 238
$ws.Range("Q160").Value = 2
$ws.Range("R160").Value = "Hortaliza"

# New "Segunda" quality observation
$ws.Range("A161").Value = 8
$ws.Range("B161").Value = "Terminal La Palmera de La Serena"
$ws.Range("C161").Value = "Coquimbo"
$ws.Range("D161").Value = 44554
$ws.Range("E161").Value = 4
$ws.Range("F161").Value = 100112009
$ws.Range("G161").Value = "Acelga"
$ws.Range("H161").Value = "Sin especificar"
$ws.Range("I161").Value = "Segunda"
$ws.Range("J161").Value = 1460
$ws.Range("K161").Value = 350
$ws.Range("L161").Value = 400
$ws.Range("M161").Value = 375
$ws.Range("N161").Value = "$/atado 1,5 a 2 kilos"
$ws.Range("O161").Value = "Provincia del Elquí"
$ws.Range("P161").Value = 188
$ws.Range("Q161").Value = 2
$ws.Range("R161").Value = "Hortaliza"
